$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 202.3077
$ws.Range("I12").Value = 150
$ws.Range("J12").Value = 286
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 286
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = -626

$ws.Range("H31").Value = 1286.4
$ws.Range("I31").Value = 363
$ws.Range("J31").Value = 4980
$ws.Range("K31").Value = 1089
$ws.Range("L31").Value = 14940
$ws.Range("M31").Value = -859
$ws.Range("N31").Value = -15400

$ws.Range("H45").Value = 4117.25
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 4117.25
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 12351.75
$ws.Range("N45").Value = -12735.75

$ws.Range("H46").Value = 3015.3157
$ws.Range("I46").Value = 1258.2
$ws.Range("J46").Value = 3642.8572
$ws.Range("K46").Value = 3774.6
$ws.Range("L46").Value = 10928.5716
$ws.Range("M46").Value = -3655.6
$ws.Range("N46").Value = -11166.5716

$ws.Range("H49").Value = 3257.7
$ws.Range("I49").Value = 859
$ws.Range("J49").Value = 4285.7144
$ws.Range("K49").Value = 2577
$ws.Range("L49").Value = 12857.1432
$ws.Range("M49").Value = -2441
$ws.Range("N49").Value = -13129.1432

$ws.Range("H59").Value = 4539.25
$ws.Range("I59").Value = 517
$ws.Range("J59").Value = 5113.857
$ws.Range("K59").Value = 1551
$ws.Range("L59").Value = 15341.571
$ws.Range("M59").Value = -994
$ws.Range("N59").Value = -16455.571

$ws.Range("H60").Value = 3015.3157
$ws.Range("I60").Value = 1258.2
$ws.Range("J60").Value = 3642.8572
$ws.Range("K60").Value = 3774.6
$ws.Range("L60").Value = 10928.5716
$ws.Range("M60").Value = -3290.6
$ws.Range("N60").Value = -11896.5716

$ws.Range("H96").Value = 325.33334
$ws.Range("I96").Value = 226
$ws.Range("J96").Value = 375
$ws.Range("K96").Value = 678
$ws.Range("L96").Value = 1125
$ws.Range("M96").Value = 695
$ws.Range("N96").Value = -3871

$ws.Range("H113").Value = 38465790
$ws.Range("I113").Value = 83335384
$ws.Range("J113").Value = 6141
$ws.Range("K113").Value = 83335384
$ws.Range("L113").Value = 6141
$ws.Range("M113").Value = -83332130
$ws.Range("N113").Value = -12649

$ws.Range("H129").Value = 1198.1666
$ws.Range("I129").Value = 519.5454999999999
$ws.Range("J129").Value = 1496.76
$ws.Range("K129").Value = 1558.6365
$ws.Range("L129").Value = 4490.28
$ws.Range("M129").Value = 3441.3635
$ws.Range("N129").Value = -14490.28

$ws.Range("H131").Value = 2110.7932
$ws.Range("I131").Value = 694.2
$ws.Range("J131").Value = 3628.5715
$ws.Range("K131").Value = 2082.6
$ws.Range("L131").Value = 10885.7145
$ws.Range("M131").Value = 2957.4
$ws.Range("N131").Value = -20965.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 37896.332
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 37896.332
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 37896.332
$ws.Range("N44").Value = -38872.332

$ws.Range("H55").Value = 23982.572
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 23982.572
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 23982.572
$ws.Range("N55").Value = -24612.572

$ws.Range("H74").Value = 903.2105
$ws.Range("I74").Value = 971.93335
$ws.Range("J74").Value = 645.5
$ws.Range("K74").Value = 971.93335
$ws.Range("L74").Value = 645.5
$ws.Range("M74").Value = -97.93335000000002
$ws.Range("N74").Value = -2393.5

$ws.Range("H77").Value = 903.2105
$ws.Range("I77").Value = 971.93335
$ws.Range("J77").Value = 645.5
$ws.Range("K77").Value = 4859.66675
$ws.Range("L77").Value = 3227.5
$ws.Range("M77").Value = -491.6667500000003
$ws.Range("N77").Value = -11963.5

$ws.Range("H80").Value = 26886.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 26886.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 26886.75
$ws.Range("N80").Value = -28882.75

$ws.Range("H83").Value = 26886.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 26886.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 80660.25
$ws.Range("N83").Value = -90644.25

$ws.Range("H97").Value = 3680
$ws.Range("I97").Value = 2673.3333
$ws.Range("J97").Value = 6700
$ws.Range("K97").Value = 2673.3333
$ws.Range("L97").Value = 6700
$ws.Range("M97").Value = -2177.3333
$ws.Range("N97").Value = -7692

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34782
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 34782
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 34782
$ws.Range("N35").Value = -35402

$ws.Range("H82").Value = 69855.63
$ws.Range("I82").Value = 205504.67
$ws.Range("J82").Value = 18987.25
$ws.Range("K82").Value = 205504.67
$ws.Range("L82").Value = 18987.25
$ws.Range("M82").Value = -205121.67
$ws.Range("N82").Value = -19753.25

$ws.Range("H85").Value = 69855.63
$ws.Range("I85").Value = 205504.67
$ws.Range("J85").Value = 18987.25
$ws.Range("K85").Value = 205504.67
$ws.Range("L85").Value = 18987.25
$ws.Range("M85").Value = -204178.67
$ws.Range("N85").Value = -21639.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 18182.5
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 18182.5
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 18182.5
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -19038.5

$ws.Range("H68").Value = 35295
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 35295
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 35295
$ws.Range("N68").Value = -36793

$ws.Range("H71").Value = 35295
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 35295
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 105885
$ws.Range("N71").Value = -113373

$ws.Range("H109").Value = 20674
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 20674
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 20674
$ws.Range("N109").Value = -22754

$ws.Range("H132").Value = 4496.1665
$ws.Range("I132").Value = 4001.7144
$ws.Range("J132").Value = 5188.4
$ws.Range("K132").Value = 12005.1432
$ws.Range("L132").Value = 15565.2
$ws.Range("M132").Value = -9475.143199999999
$ws.Range("N132").Value = -20625.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 38531.848
$ws.Range("I2").Value = 142898
$ws.Range("J2").Value = 81.1579
$ws.Range("K2").Value = 857388
$ws.Range("L2").Value = 486.9474
$ws.Range("M2").Value = -857275
$ws.Range("N2").Value = -712.9474

$ws.Range("H75").Value = 2315.3845
$ws.Range("I75").Value = 350
$ws.Range("J75").Value = 2672.7273
$ws.Range("K75").Value = 1050
$ws.Range("L75").Value = 8018.1819
$ws.Range("M75").Value = -52
$ws.Range("N75").Value = -10014.1819

$ws.Range("H78").Value = 2315.3845
$ws.Range("I78").Value = 350
$ws.Range("J78").Value = 2672.7273
$ws.Range("K78").Value = 3150
$ws.Range("L78").Value = 24054.5457
$ws.Range("M78").Value = 1842
$ws.Range("N78").Value = -34038.5457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19227.25
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 19227.25
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 19227.25
$ws.Range("N57").Value = -20867.25

$ws.Range("H63").Value = 20000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 20000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372

$ws.Range("H66").Value = 20000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 20000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 20867.727
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 20867.727
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 20867.727
$ws.Range("N63").Value = -22365.727

$ws.Range("H66").Value = 20867.727
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 20867.727
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 62603.181
$ws.Range("N66").Value = -70091.181

$ws.Range("H69").Value = 25000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 25000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26622

$ws.Range("H72").Value = 25000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 25000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -83112

$ws.Range("H132").Value = 2451.8572
$ws.Range("I132").Value = 2207.2727
$ws.Range("J132").Value = 2720.9
$ws.Range("K132").Value = 6621.8181
$ws.Range("L132").Value = 8162.700000000001
$ws.Range("M132").Value = -4091.8181
$ws.Range("N132").Value = -13222.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 27478.334
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 27478.334
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 27478.334
$ws.Range("N64").Value = -27974.334

$ws.Range("H67").Value = 27478.334
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 27478.334
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 27478.334
$ws.Range("N67").Value = -29194.334

$ws.Range("H100").Value = 833932.4399999999
$ws.Range("I100").Value = 631.7778
$ws.Range("J100").Value = 3333834.2
$ws.Range("K100").Value = 1263.5556
$ws.Range("L100").Value = 6667668.4
$ws.Range("M100").Value = -722.5555999999999
$ws.Range("N100").Value = -6668750.4

$ws.Range("H109").Value = 29438.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 29438.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 29438.5
$ws.Range("N109").Value = -32212.5

Write-Output "Applied profit sheet updates"
